$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column (B) stays as literal text instead of being
# auto-converted to a date serial when we assign "yyyy-mm-dd" strings.
$ws.Range("B2:B18").NumberFormat = "@"

# ---- Simple single-cell odds corrections (rows 2,3,4,5,8,9,13) ----
$ws.Range("X2").Value = 980
$ws.Range("R3").Value = 1.63
$ws.Range("P4").Value = 2.12
$ws.Range("U5").Value = 2.62
$ws.Range("AC5").Value = 8.8
$ws.Range("AD5").Value = 11
$ws.Range("G8").Value = 3.05
$ws.Range("W8").Value = 1.48
$ws.Range("J9").Value = 2.86
$ws.Range("K9").Value = 3.1
$ws.Range("N9").Value = 2.34
$ws.Range("O9").Value = 1.66
$ws.Range("P9").Value = 1.44
$ws.Range("T9").Value = 2.3
$ws.Range("X9").Value = 7.2
$ws.Range("AN9").Value = 55
$ws.Range("T13").Value = 2.32

# ---- Rows 10 and 11 swap content (Paraguayan Primera Division <-> Colombian Primera A)
# with a couple of odds tweaks in the relocated Colombian Primera A row ----
$block1 = New-Object 'object[,]' 2,41
$block1[0,0] = 'Paraguayan Primera Division'
$block1[0,1] = '2026-01-27'
$block1[0,2] = '18:00:00'
$block1[0,3] = 'Sportivo Luqueno'
$block1[0,4] = 'Nacional (Par)'
$block1[0,5] = 0
$block1[0,6] = 0
$block1[0,7] = 0
$block1[0,8] = 0
$block1[0,9] = 0
$block1[0,10] = 0
$block1[0,11] = 1.01
$block1[0,12] = 1.01
$block1[0,13] = 1.02
$block1[0,14] = 1.42
$block1[0,15] = 1.07
$block1[0,16] = 1.01
$block1[0,17] = 1.08
$block1[0,18] = 1.01
$block1[0,19] = 1.01
$block1[0,20] = 1.01
$block1[0,21] = 1.01
$block1[0,22] = 1.01
$block1[0,23] = 1000
$block1[0,24] = 1000
$block1[0,25] = 1000
$block1[0,26] = 1000
$block1[0,27] = 1000
$block1[0,28] = 1000
$block1[0,29] = 1000
$block1[0,30] = 1000
$block1[0,31] = 1000
$block1[0,32] = 1000
$block1[0,33] = 1000
$block1[0,34] = 1000
$block1[0,35] = 1000
$block1[0,36] = 1000
$block1[0,37] = 1000
$block1[0,38] = 1000
$block1[0,39] = 1000
$block1[0,40] = 1000
$block1[1,0] = 'Colombian Primera A'
$block1[1,1] = '2026-01-27'
$block1[1,2] = '18:00:00'
$block1[1,3] = 'Cucuta Deportivo'
$block1[1,4] = 'Atletico Bucaramanga'
$block1[1,5] = 2.68
$block1[1,6] = 3.6
$block1[1,7] = 2.42
$block1[1,8] = 3.25
$block1[1,9] = 2.78
$block1[1,10] = 3.6
$block1[1,11] = 1.01
$block1[1,12] = 1.01
$block1[1,13] = 2.5
$block1[1,14] = 1.53
$block1[1,15] = 1.5
$block1[1,16] = 2.42
$block1[1,17] = 1.14
$block1[1,18] = 4.5
$block1[1,19] = 1.78
$block1[1,20] = 1.56
$block1[1,21] = 1.44
$block1[1,22] = 1.38
$block1[1,23] = 11
$block1[1,24] = 970
$block1[1,25] = 970
$block1[1,26] = 1000
$block1[1,27] = 970
$block1[1,28] = 970
$block1[1,29] = 970
$block1[1,30] = 1000
$block1[1,31] = 970
$block1[1,32] = 970
$block1[1,33] = 1000
$block1[1,34] = 1000
$block1[1,35] = 1000
$block1[1,36] = 1000
$block1[1,37] = 1000
$block1[1,38] = 1000
$block1[1,39] = 1000
$block1[1,40] = 1000
$ws.Range("A10:AO11").Value = $block1

# ---- Row 14 (Fortaleza FC / Llaneros FC) shifts down to row 15 (with small odds
# tweaks), a new Paraguayan Primera Division match is inserted as row 14, and three
# brand new Argentinian Primera Division / Colombian Primera A matches are appended
# as rows 16-18 ----
$block2 = New-Object 'object[,]' 5,41
$block2[0,0] = 'Paraguayan Primera Division'
$block2[0,1] = '2026-01-27'
$block2[0,2] = '20:15:00'
$block2[0,3] = 'Sportivo San Lorenzo'
$block2[0,4] = 'Cerro Porteno'
$block2[0,5] = 0
$block2[0,6] = 0
$block2[0,7] = 0
$block2[0,8] = 0
$block2[0,9] = 0
$block2[0,10] = 0
$block2[0,11] = 1.01
$block2[0,12] = 1.01
$block2[0,13] = 1.02
$block2[0,14] = 1.28
$block2[0,15] = 1.07
$block2[0,16] = 1.01
$block2[0,17] = 1.08
$block2[0,18] = 1.01
$block2[0,19] = 1.01
$block2[0,20] = 1.01
$block2[0,21] = 1.01
$block2[0,22] = 1.01
$block2[0,23] = 1000
$block2[0,24] = 1000
$block2[0,25] = 1000
$block2[0,26] = 1000
$block2[0,27] = 1000
$block2[0,28] = 1000
$block2[0,29] = 1000
$block2[0,30] = 1000
$block2[0,31] = 1000
$block2[0,32] = 1000
$block2[0,33] = 1000
$block2[0,34] = 1000
$block2[0,35] = 1000
$block2[0,36] = 1000
$block2[0,37] = 1000
$block2[0,38] = 1000
$block2[0,39] = 1000
$block2[0,40] = 1000
$block2[1,0] = 'Colombian Primera A'
$block2[1,1] = '2026-01-27'
$block2[1,2] = '20:20:00'
$block2[1,3] = 'Fortaleza FC'
$block2[1,4] = 'Llaneros FC'
$block2[1,5] = 1.93
$block2[1,6] = 2.12
$block2[1,7] = 4.5
$block2[1,8] = 5.9
$block2[1,9] = 3.1
$block2[1,10] = 3.6
$block2[1,11] = 1.01
$block2[1,12] = 1.08
$block2[1,13] = 2.9
$block2[1,14] = 1.43
$block2[1,15] = 1.59
$block2[1,16] = 2.18
$block2[1,17] = 1.21
$block2[1,18] = 3.9
$block2[1,19] = 1.01
$block2[1,20] = 1.76
$block2[1,21] = 1.24
$block2[1,22] = 1.9
$block2[1,23] = 970
$block2[1,24] = 19.5
$block2[1,25] = 48
$block2[1,26] = 1000
$block2[1,27] = 10
$block2[1,28] = 10
$block2[1,29] = 28
$block2[1,30] = 100
$block2[1,31] = 16
$block2[1,32] = 15
$block2[1,33] = 30
$block2[1,34] = 1000
$block2[1,35] = 36
$block2[1,36] = 36
$block2[1,37] = 70
$block2[1,38] = 1000
$block2[1,39] = 1000
$block2[1,40] = 1000
$block2[2,0] = 'Argentinian Primera Division'
$block2[2,1] = '2026-01-27'
$block2[2,2] = '22:15:00'
$block2[2,3] = 'Newells'
$block2[2,4] = 'CA Independiente'
$block2[2,5] = 3.55
$block2[2,6] = 5.1
$block2[2,7] = 2.14
$block2[2,8] = 2.74
$block2[2,9] = 2.28
$block2[2,10] = 4.1
$block2[2,11] = 1.01
$block2[2,12] = 1.01
$block2[2,13] = 1.45
$block2[2,14] = 1.02
$block2[2,15] = 1.45
$block2[2,16] = 2.88
$block2[2,17] = 1.13
$block2[2,18] = 5
$block2[2,19] = 2.06
$block2[2,20] = 1.54
$block2[2,21] = 1.57
$block2[2,22] = 1.25
$block2[2,23] = 9.4
$block2[2,24] = 9.8
$block2[2,25] = 20
$block2[2,26] = 55
$block2[2,27] = 14.5
$block2[2,28] = 10
$block2[2,29] = 18.5
$block2[2,30] = 55
$block2[2,31] = 970
$block2[2,32] = 970
$block2[2,33] = 970
$block2[2,34] = 100
$block2[2,35] = 1000
$block2[2,36] = 100
$block2[2,37] = 1000
$block2[2,38] = 1000
$block2[2,39] = 1000
$block2[2,40] = 1000
$block2[3,0] = 'Argentinian Primera Division'
$block2[3,1] = '2026-01-27'
$block2[3,2] = '22:15:00'
$block2[3,3] = 'Atl Tucuman'
$block2[3,4] = 'Central Cordoba (SdE)'
$block2[3,5] = 2.18
$block2[3,6] = 2.8
$block2[3,7] = 3.6
$block2[3,8] = 5.6
$block2[3,9] = 2.16
$block2[3,10] = 3.7
$block2[3,11] = 1.01
$block2[3,12] = 1.01
$block2[3,13] = 1.13
$block2[3,14] = 1.02
$block2[3,15] = 1.13
$block2[3,16] = 1.6
$block2[3,17] = 1.13
$block2[3,18] = 4.8
$block2[3,19] = 1.84
$block2[3,20] = 1.51
$block2[3,21] = 1.22
$block2[3,22] = 1.55
$block2[3,23] = 9.6
$block2[3,24] = 13.5
$block2[3,25] = 36
$block2[3,26] = 1000
$block2[3,27] = 9.6
$block2[3,28] = 9.8
$block2[3,29] = 25
$block2[3,30] = 95
$block2[3,31] = 19
$block2[3,32] = 17.5
$block2[3,33] = 34
$block2[3,34] = 1000
$block2[3,35] = 50
$block2[3,36] = 50
$block2[3,37] = 95
$block2[3,38] = 1000
$block2[3,39] = 1000
$block2[3,40] = 1000
$block2[4,0] = 'Colombian Primera A'
$block2[4,1] = '2026-01-27'
$block2[4,2] = '22:30:00'
$block2[4,3] = 'Ind Medellin'
$block2[4,4] = 'Tolima'
$block2[4,5] = 2.06
$block2[4,6] = 2.2
$block2[4,7] = 4.2
$block2[4,8] = 4.9
$block2[4,9] = 3.05
$block2[4,10] = 3.4
$block2[4,11] = 1.01
$block2[4,12] = 1.01
$block2[4,13] = 2.8
$block2[4,14] = 1.45
$block2[4,15] = 1.59
$block2[4,16] = 2.4
$block2[4,17] = 1.19
$block2[4,18] = 4
$block2[4,19] = 1.01
$block2[4,20] = 1.76
$block2[4,21] = 1.25
$block2[4,22] = 1.83
$block2[4,23] = 14.5
$block2[4,24] = 17.5
$block2[4,25] = 44
$block2[4,26] = 1000
$block2[4,27] = 10
$block2[4,28] = 10
$block2[4,29] = 28
$block2[4,30] = 100
$block2[4,31] = 17
$block2[4,32] = 16
$block2[4,33] = 32
$block2[4,34] = 1000
$block2[4,35] = 40
$block2[4,36] = 40
$block2[4,37] = 75
$block2[4,38] = 1000
$block2[4,39] = 1000
$block2[4,40] = 1000
$ws.Range("A14:AO18").Value = $block2
